$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1867704280155642
$ws.Range("C2").Value = 0.5603112840466926
$ws.Range("J2").Value = 0.01167315175097276
$ws.Range("P2").Value = 0.1478599221789883
$ws.Range("S2").Value = 0.0933852140077821
# Row 3
$ws.Range("B3").Value = 0.0131578947368421
$ws.Range("C3").Value = 0.03947368421052631
$ws.Range("J3").Value = 0.03947368421052631
$ws.Range("P3").Value = 0.6842105263157895
$ws.Range("S3").Value = 0.2236842105263158
# Row 4
$ws.Range("P4").Value = 0.7777777777777778
$ws.Range("S4").Value = 0.2222222222222222
# Row 6
$ws.Range("B6").Value = 0.06829268292682927
$ws.Range("F6").Value = 0.06829268292682927
$ws.Range("J6").Value = 0.2585365853658537
$ws.Range("O6").Value = 0.00975609756097561
$ws.Range("Q6").Value = 0.1121951219512195
$ws.Range("R6").Value = 0.07317073170731707
$ws.Range("S6").Value = 0.4097560975609756
# Row 7
$ws.Range("B7").Value = 0.07065217391304347
$ws.Range("D7").Value = 0.0108695652173913
$ws.Range("E7").Value = 0.0108695652173913
$ws.Range("F7").Value = 0.09782608695652174
$ws.Range("J7").Value = 0.1032608695652174
$ws.Range("O7").Value = 0.03260869565217391
$ws.Range("Q7").Value = 0.1304347826086956
$ws.Range("R7").Value = 0.09239130434782608
$ws.Range("S7").Value = 0.4510869565217391
# Row 8
$ws.Range("B8").Value = 0.1135371179039301
$ws.Range("D8").Value = 0.01746724890829694
$ws.Range("F8").Value = 0.05895196506550218
$ws.Range("J8").Value = 0.1091703056768559
$ws.Range("O8").Value = 0.04366812227074236
$ws.Range("Q8").Value = 0.1069868995633188
$ws.Range("R8").Value = 0.1004366812227074
$ws.Range("S8").Value = 0.4497816593886463
# Row 9
$ws.Range("B9").Value = 0.1111111111111111
$ws.Range("J9").Value = 0.08496732026143791
$ws.Range("O9").Value = 0.0261437908496732
$ws.Range("Q9").Value = 0.1503267973856209
$ws.Range("S9").Value = 0.4509803921568628
# Row 10
$ws.Range("B10").Value = 0.1055226824457594
$ws.Range("D10").Value = 0.009861932938856016
$ws.Range("E10").Value = 0.0009861932938856016
$ws.Range("F10").Value = 0.07988165680473373
$ws.Range("J10").Value = 0.141025641025641
$ws.Range("O10").Value = 0.02465483234714004
$ws.Range("Q10").Value = 0.1627218934911243
$ws.Range("R10").Value = 0.08382642998027613
$ws.Range("S10").Value = 0.3915187376725838
# Row 11
$ws.Range("G11").Value = 0.1418685121107267
$ws.Range("J11").Value = 0.09342560553633218
$ws.Range("K11").Value = 0.2076124567474048
$ws.Range("L11").Value = 0.532871972318339
$ws.Range("S11").Value = 0.02422145328719723
# Row 12
$ws.Range("G12").Value = 0.6727272727272727
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.01818181818181818
$ws.Range("L12").Value = 0.05454545454545454
$ws.Range("S12").Value = 0.05454545454545454
# Row 13
$ws.Range("G13").Value = 0.7592592592592593
$ws.Range("J13").Value = 0.1851851851851852
# Row 15
$ws.Range("F15").Value = 0.005050505050505051
$ws.Range("H15").Value = 0.1515151515151515
$ws.Range("I15").Value = 0.06565656565656566
$ws.Range("J15").Value = 0.2474747474747475
$ws.Range("K15").Value = 0.06060606060606061
$ws.Range("M15").Value = 0.005050505050505051
$ws.Range("O15").Value = 0.09090909090909091
$ws.Range("S15").Value = 0.3737373737373738
# Row 16
$ws.Range("F16").Value = 0.01282051282051282
$ws.Range("H16").Value = 0.1474358974358974
$ws.Range("I16").Value = 0.07692307692307693
$ws.Range("J16").Value = 0.391025641025641
$ws.Range("K16").Value = 0.1282051282051282
$ws.Range("M16").Value = 0.04487179487179487
$ws.Range("O16").Value = 0.03205128205128205
$ws.Range("S16").Value = 0.1666666666666667
# Row 17
$ws.Range("F17").Value = 0.01428571428571429
$ws.Range("H17").Value = 0.2035714285714286
$ws.Range("I17").Value = 0.06071428571428571
$ws.Range("J17").Value = 0.4035714285714286
$ws.Range("K17").Value = 0.1321428571428571
$ws.Range("M17").Value = 0.02142857142857143
$ws.Range("O17").Value = 0.075
$ws.Range("S17").Value = 0.08928571428571429
# Row 18
$ws.Range("F18").Value = 0.02777777777777778
$ws.Range("H18").Value = 0.2055555555555555
$ws.Range("I18").Value = 0.06666666666666667
$ws.Range("J18").Value = 0.3888888888888889
$ws.Range("K18").Value = 0.1166666666666667
$ws.Range("M18").Value = 0.01666666666666667
$ws.Range("O18").Value = 0.05
$ws.Range("S18").Value = 0.1277777777777778
# Row 19
$ws.Range("F19").Value = 0.0183180682764363
$ws.Range("H19").Value = 0.2547876769358868
$ws.Range("I19").Value = 0.08159866777685262
$ws.Range("J19").Value = 0.324729392173189
$ws.Range("K19").Value = 0.1115736885928393
$ws.Range("M19").Value = 0.03330557868442964
$ws.Range("N19").Value = 0.001665278934221482
$ws.Range("O19").Value = 0.05328892589508743
$ws.Range("S19").Value = 0.1207327227310575
